$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "submitted" column (C) to TRUE for rows 3 through 70 (row 2 already TRUE).
for ($r = 3; $r -le 70; $r++) {
    $ws.Cells.Item($r, 3).Value = $true
}

# Update the selection/view: active cell C3, single-cell selection, no frozen/scrolled topLeftCell.
$ws.Range("C3").Select()
